$d = $word.ActiveDocument

# Locate the empty "List Paragraph" item that sits right before the
# page-break paragraph at the end of the Case 1 bullet list (the one
# that carries the now-invisible _GoBack bookmark). We find it by
# scanning for the run of text "Misschien alternatieve van Markdown"
# (the bullet right before it) and taking the following paragraph.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "Misschien alternatieve van Markdown") {
        $target = $i + 1
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate target paragraph"
}

# Step 1: give the (currently empty) bullet its text.
$p1 = $d.Paragraphs.Item($target)
$r1 = $p1.Range
$ip1 = $d.Range($r1.End - 1, $r1.End - 1)
$ip1.InsertAfter("Of we een login nodig hebben")

# Step 2: split off a brand new bullet paragraph after it.
$p1b = $d.Paragraphs.Item($target)
$r1b = $p1b.Range
$ip2 = $d.Range($r1b.End - 1, $r1b.End - 1)
$ip2.InsertParagraphAfter()

# Step 3: the newly created (still empty) bullet paragraph now holds
# the bookmark that used to live on the original paragraph; give it
# the "Planning" text.
$p2 = $d.Paragraphs.Item($target + 1)
$r2 = $p2.Range
$ip3 = $d.Range($r2.End - 1, $r2.End - 1)
$ip3.InsertAfter("Planning")
